$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-22 down to 11-23.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with the new weekly data point,
# matching the structure/format of the surrounding rows.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44883
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112026
$ws.Range("G10").Value = "Haba"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 550
$ws.Range("L10").Value = 600
$ws.Range("M10").Value = 575
$ws.Range("N10").Value = "$/kilo"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 575
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
